$p = $ppt.ActivePresentation

# 1. Remove the first (empty) title slide - only the diagrams slide remains.
$p.Slides.Item(1).Delete()

# 2. Bump the cached "today" date shown by the auto date placeholder
#    (master + every layout) from 12/7/2014 to 12/8/2014.
function Set-DatePlaceholderText($shapes, $newText) {
  for ($si = 1; $si -le $shapes.Count; $si++) {
    $sh = $shapes.Item($si)
    $isDatePh = $false
    try {
      if ($sh.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
    } catch {}
    if ($isDatePh) {
      $sh.TextFrame.TextRange.Text = $newText
    }
  }
}

$m = $p.SlideMaster
Set-DatePlaceholderText $m.Shapes "12/8/2014"
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
  $lay = $m.CustomLayouts.Item($li)
  Set-DatePlaceholderText $lay.Shapes "12/8/2014"
}
